$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column (BF) was populated using the wrong date format
# (the raw "5-14-2012-13" string mashed together from the season label and
# the game date). Correct it to an ISO-formatted date string "2013-05-14"
# for every data row (the data was one day off due to how NBA stats were
# displayed/scraped).
$rng = $ws.Range("BF2:BF31")

# Force a text number format first so Excel keeps the corrected value as a
# literal string instead of re-interpreting "2013-05-14" as a date serial.
$rng.NumberFormat = "@"
$rng.Value = "2013-05-14"
